$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price column to Text format so purely-numeric-looking
# values (e.g. "0.998") are stored as strings, matching the source data,
# then restore the default (unstyled) cell style once values are set.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Cells.Item(2, 4).Value = '70.191.73'
$ws.Cells.Item(2, 5).Value = '  +0.49%  '
$ws.Cells.Item(3, 4).Value = '3.605.37'
$ws.Cells.Item(3, 5).Value = '  +2.43%  '
$ws.Cells.Item(4, 4).Value = '0.998'
$ws.Cells.Item(4, 5).Value = '  -0.04%  '
$ws.Cells.Item(5, 4).Value = '604.59'
$ws.Cells.Item(5, 5).Value = '  +0.55%  '
$ws.Cells.Item(6, 4).Value = '196.58'
$ws.Cells.Item(6, 5).Value = '  +0.59%  '
$ws.Cells.Item(7, 4).Value = '0.626'
$ws.Cells.Item(7, 5).Value = '  +0.20%  '
$ws.Cells.Item(8, 5).Value = '  +0.06%  '
$ws.Cells.Item(9, 5).Value = '  -1.82%  '
$ws.Cells.Item(10, 5).Value = '  -0.82%  '
$ws.Cells.Item(11, 4).Value = '53.82'
$ws.Cells.Item(11, 5).Value = '  -0.40%  '
$ws.Cells.Item(12, 5).Value = '  +0.73%  '
$ws.Cells.Item(13, 4).Value = '9.59'
$ws.Cells.Item(13, 5).Value = '  +0.65%  '
$ws.Cells.Item(14, 4).Value = '4.178.88'
$ws.Cells.Item(14, 5).Value = '  +2.56%  '
$ws.Cells.Item(15, 4).Value = '13.09'
$ws.Cells.Item(15, 5).Value = '  +3.69%  '
$ws.Cells.Item(16, 4).Value = '596.37'
$ws.Cells.Item(16, 5).Value = '  -1.12%  '
$ws.Cells.Item(17, 4).Value = '70.330.46'
$ws.Cells.Item(17, 5).Value = '  +0.41%  '
$ws.Cells.Item(18, 4).Value = '19.11'
$ws.Cells.Item(18, 5).Value = '  -0.08%  '
$ws.Cells.Item(19, 4).Value = '3.601.46'
$ws.Cells.Item(19, 5).Value = '  +2.36%  '
$ws.Cells.Item(20, 5).Value = '  +1.39%  '
$ws.Cells.Item(21, 4).Value = '0.996'
$ws.Cells.Item(21, 5).Value = '  +0.17%  '
$ws.Cells.Item(22, 4).Value = '17.77'
$ws.Cells.Item(22, 5).Value = '  -2.83%  '
$ws.Cells.Item(23, 4).Value = '5.17'
$ws.Cells.Item(23, 5).Value = '  -1.08%  '
$ws.Cells.Item(24, 4).Value = '102.02'
$ws.Cells.Item(24, 5).Value = '  -1.72%  '
$ws.Cells.Item(25, 4).Value = '4.61'
$ws.Cells.Item(25, 5).Value = '  +0.22%  '
$ws.Cells.Item(26, 4).Value = '3.02'
$ws.Cells.Item(26, 5).Value = '  -1.51%  '
$ws.Cells.Item(27, 4).Value = '10.74'
$ws.Cells.Item(27, 5).Value = '  -1.60%  '
$ws.Cells.Item(28, 4).Value = '9.59'
$ws.Cells.Item(28, 5).Value = '  -0.78%  '
$ws.Cells.Item(29, 4).Value = '33.79'
$ws.Cells.Item(29, 5).Value = '  +0.72%  '
$ws.Cells.Item(30, 4).Value = '4.73'
$ws.Cells.Item(30, 5).Value = '  +5.37%  '
$ws.Cells.Item(31, 5).Value = '  +0.70%  '
$ws.Cells.Item(32, 5).Value = '  -3.35%  '
$ws.Cells.Item(33, 5).Value = '  +0.91%  '
$ws.Cells.Item(34, 4).Value = '63.30'
$ws.Cells.Item(34, 5).Value = '  +0.14%  '
$ws.Cells.Item(35, 4).Value = '0.0₃0886'
$ws.Cells.Item(35, 5).Value = '  +7.51%  '
$ws.Cells.Item(36, 4).Value = '3.904.14'
$ws.Cells.Item(36, 5).Value = '  +3.68%  '
$ws.Cells.Item(37, 5).Value = '  +0.70%  '
$ws.Cells.Item(38, 5).Value = '  -0.01%  '
$ws.Cells.Item(39, 4).Value = '517.48'
$ws.Cells.Item(39, 5).Value = '  +6.03%  '
$ws.Cells.Item(40, 4).Value = '36.87'
$ws.Cells.Item(40, 5).Value = '  +0.27%  '
$ws.Cells.Item(41, 5).Value = '  -0.92%  '
$ws.Cells.Item(42, 5).Value = '  -1.93%  '
$ws.Cells.Item(43, 5).Value = '  -2.03%  '
$ws.Cells.Item(44, 5).Value = '  -0.41%  '
$ws.Cells.Item(45, 4).Value = '3.43'
$ws.Cells.Item(45, 5).Value = '  +3.25%  '
$ws.Cells.Item(46, 5).Value = '  +1.31%  '
$ws.Cells.Item(47, 5).Value = '  -0.06%  '
$ws.Cells.Item(48, 4).Value = '8.63'
$ws.Cells.Item(48, 5).Value = '  -0.20%  '
$ws.Cells.Item(49, 5).Value = '  -0.36%  '
$ws.Cells.Item(50, 4).Value = '0.000250'
$ws.Cells.Item(50, 5).Value = '  +2.95%  '
$ws.Cells.Item(51, 5).Value = '  -1.48%  '

# Restore the original (default) style on the Price column so no
# extraneous style index is left attached to the cells.
$ws.Range("D2:D51").Style = "Normal"
